# Update the time schedule
# Expand the Project Milestones table with new controller-design milestones
# and refresh the dates/positions of the existing ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Make room for the 5 new milestone rows inside the table.
#    (Row numbers below refer to the CURRENT sheet at the moment each
#    Insert() runs - i.e. they already account for the shifting caused
#    by the previous inserts.)
# ------------------------------------------------------------------

# Row 19 used to be "MCP" - it will be overwritten with "MVC" below, and
# we need one extra row right after it for "Simulation".
$ws.Rows.Item(20).Insert()

# Two extra rows are needed before "STOK" (currently row 22) for the
# "Controller design - PI" / "Controller design - LQR" milestones.
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

# Two extra rows are needed right after "STOK" (now row 24) for the
# "Controller design - MPC" / "Controller Test" milestones.
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(25).Insert()

# ------------------------------------------------------------------
# 2. Write the refreshed milestone data (Date / Milestone / Position).
# ------------------------------------------------------------------

$milestones = @(
    @{ Row = 17; Date = 44470; Name = "Project Start";            Pos = 10 },
    @{ Row = 18; Date = 44489; Name = "Model ready";               Pos = 20 },
    @{ Row = 19; Date = 44495; Name = "MVC";                       Pos = -10 },
    @{ Row = 20; Date = 44496; Name = "Simulation";                Pos = 20 },
    @{ Row = 21; Date = 44501; Name = "Measure in LAB";            Pos = 10 },
    @{ Row = 22; Date = 44505; Name = "Controller design - PI";    Pos = 20 },
    @{ Row = 23; Date = 44512; Name = "Controller design - LQR";   Pos = 10 },
    @{ Row = 24; Date = 44516; Name = "STOK";                      Pos = -10 },
    @{ Row = 25; Date = 44518; Name = "Controller design - MPC";   Pos = 20 },
    @{ Row = 26; Date = 44525; Name = "Controller Test";           Pos = 10 },
    @{ Row = 27; Date = 44526; Name = "Artikel+Abstract";          Pos = 20 },
    @{ Row = 28; Date = 44529; Name = "DRTS";                      Pos = -10 },
    @{ Row = 29; Date = 44538; Name = "Review(s)";                 Pos = 20 },
    @{ Row = 30; Date = 44545; Name = "Project done";              Pos = 20 },
    @{ Row = 31; Date = 44547; Name = "Project delivery";          Pos = 10 },
    @{ Row = 32; Date = 44551; Name = "SEMCON";                    Pos = 20 },
    @{ Row = 33; Date = 44562; Name = "Project End";               Pos = 10 }
)

foreach ($m in $milestones) {
    $r = $m.Row
    $ws.Range("B$r").Value = $m.Date
    $ws.Range("C$r").Value = $m.Name
    $ws.Range("E$r").Value = $m.Pos
    $ws.Range("F$r").Formula = "=0"
}

# ------------------------------------------------------------------
# 3. Resize the ProjectDetails table to cover the new data + trailing
#    blank row (B16:F34).
# ------------------------------------------------------------------

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B16:F34"))

# ------------------------------------------------------------------
# 4. Keep the milestone "tip" box (merged H17:L..) in sync - it now
#    spans down to row 24 to match the taller milestone list.
# ------------------------------------------------------------------

$ws.Range("H17:L24").Merge()

# ------------------------------------------------------------------
# 5. Point the chart at the new data range so the series formulas
#    reference the full, expanded table.
# ------------------------------------------------------------------

$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.SeriesCollection(1).Formula = "=SERIES('Project Timeline'!`$E`$16,'Project Timeline'!`$C`$17:`$C`$33,'Project Timeline'!`$E`$17:`$E`$33,1)"
$chart.SeriesCollection(2).Formula = "=SERIES('Project Timeline'!`$B`$16,'Project Timeline'!`$B`$17:`$B`$33,'Project Timeline'!`$F`$17:`$F`$34,2)"
